$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove three course rows that no longer appear in the updated dataset.
# Deleting from the bottom up so earlier row numbers stay valid as we go.
$ws.Rows.Item(24).Delete()  # PROYECTO DE SISTEMAS ROBUSTOS, PARALELOS Y DISTRIBUIDOS
$ws.Rows.Item(23).Delete()  # PROYECTO DE GESTION DE LA TECNOLOGIA DE INFORMACION
$ws.Rows.Item(9).Delete()   # COMPUTO FLEXIBLE (SOFTCOMPUTING)
